$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 60685
$ws.Range("E2").Value = -204
$ws.Range("F2").Value = -204
$ws.Range("G2").Value = -3562
$ws.Range("H2").Value = -2925
$ws.Range("I2").Value = -2743
$ws.Range("J2").Value = -182
$ws.Range("K2").Value = 88537
$ws.Range("L2").Value = 62460
$ws.Range("M2").Value = 26078
$ws.Range("N2").Value = 21965
$ws.Range("O2").Value = 4113
$ws.Range("P2").Value = 5562
$ws.Range("Q2").Value = 687
$ws.Range("R2").Value = -3752
$ws.Range("S2").Value = 483
$ws.Range("T2").Value = 2133
$ws.Range("U2").Value = -1446
$ws.Range("V2").Value = 51948
$ws.Range("W2").Value = -0.34
$ws.Range("X2").Value = -4.82
$ws.Range("Y2").Value = -12.25
$ws.Range("Z2").Value = -3.26
$ws.Range("AA2").Value = 239.51
$ws.Range("AB2").Value = 196.76
$ws.Range("AC2").Value = -3479
$ws.Range("AD2").Value = -1.68
$ws.Range("AE2").Value = 25080
$ws.Range("AF2").Value = 0.23
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 88824290

# Row 3
$ws.Range("D3").Value = 52663
$ws.Range("E3").Value = 1694
$ws.Range("F3").Value = 1936
$ws.Range("G3").Value = -2145
$ws.Range("H3").Value = -2244
$ws.Range("I3").Value = -2387
$ws.Range("J3").Value = 143
$ws.Range("K3").Value = 68959
$ws.Range("L3").Value = 46500
$ws.Range("M3").Value = 22459
$ws.Range("N3").Value = 20551
$ws.Range("O3").Value = 1908
$ws.Range("P3").Value = 5889
$ws.Range("Q3").Value = 9592
$ws.Range("R3").Value = 5669
$ws.Range("S3").Value = -14759
$ws.Range("T3").Value = 920
$ws.Range("U3").Value = 8673
$ws.Range("V3").Value = 37230
$ws.Range("W3").Value = 3.22
$ws.Range("X3").Value = -4.26
$ws.Range("Y3").Value = -11.23
$ws.Range("Z3").Value = -2.85
$ws.Range("AA3").Value = 207.04
$ws.Range("AB3").Value = 143.68
$ws.Range("AC3").Value = -2503
$ws.Range("AD3").Value = -2.3
$ws.Range("AE3").Value = 22004
$ws.Range("AF3").Value = 0.26
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 95358542

# Row 4
$ws.Range("D4").Value = 50066
$ws.Range("E4").Value = 2566
$ws.Range("F4").Value = 2566
$ws.Range("G4").Value = 459
$ws.Range("H4").Value = 708
$ws.Range("I4").Value = 613
$ws.Range("J4").Value = 95
$ws.Range("K4").Value = 63247
$ws.Range("L4").Value = 40353
$ws.Range("M4").Value = 22893
$ws.Range("N4").Value = 21796
$ws.Range("O4").Value = 1097
$ws.Range("P4").Value = 5892
$ws.Range("Q4").Value = 8883
$ws.Range("R4").Value = -1012
$ws.Range("S4").Value = -6795
$ws.Range("T4").Value = 1123
$ws.Range("U4").Value = 7760
$ws.Range("V4").Value = 30879
$ws.Range("W4").Value = 5.13
$ws.Range("X4").Value = 1.41
$ws.Range("Y4").Value = 2.9
$ws.Range("Z4").Value = 1.07
$ws.Range("AA4").Value = 176.27
$ws.Range("AB4").Value = 154.23
$ws.Range("AC4").Value = 643
$ws.Range("AD4").Value = 17.2
$ws.Range("AE4").Value = 23347
$ws.Range("AF4").Value = 0.47
$ws.Range("AG4").Value = 150
$ws.Range("AH4").Value = 1.36
$ws.Range("AI4").Value = 22.85
$ws.Range("AJ4").Value = 95420177

# Row 5
$ws.Range("D5").Value = 60493
$ws.Range("E5").Value = 2413
$ws.Range("F5").Value = 2413
$ws.Range("G5").Value = 377
$ws.Range("H5").Value = 48
$ws.Range("I5").Value = 45
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 60910
$ws.Range("L5").Value = 37203
$ws.Range("M5").Value = 23707
$ws.Range("N5").Value = 22540
$ws.Range("O5").Value = 1166
$ws.Range("P5").Value = 5892
$ws.Range("Q5").Value = 902
$ws.Range("R5").Value = -145
$ws.Range("S5").Value = -1831
$ws.Range("T5").Value = 596
$ws.Range("U5").Value = 306
$ws.Range("V5").Value = 28632
$ws.Range("W5").Value = 3.99
$ws.Range("X5").Value = 0.08
$ws.Range("Y5").Value = 0.2
$ws.Range("Z5").Value = 0.08
$ws.Range("AA5").Value = 156.93
$ws.Range("AB5").Value = 156.1
$ws.Range("AC5").Value = 47
$ws.Range("AD5").Value = 231.73
$ws.Range("AE5").Value = 24141
$ws.Range("AF5").Value = 0.46
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 0.91
$ws.Range("AI5").Value = 206.11
$ws.Range("AJ5").Value = 95432737

# Row 6
$ws.Range("D6").Value = 59649
$ws.Range("E6").Value = 1450
$ws.Range("F6").Value = 1450
$ws.Range("G6").Value = -3270
$ws.Range("H6").Value = -3045
$ws.Range("I6").Value = -2979
$ws.Range("K6").Value = 57302
$ws.Range("L6").Value = 36980
$ws.Range("M6").Value = 20322
$ws.Range("N6").Value = 19257
$ws.Range("P6").Value = 5892
$ws.Range("Q6").Value = 2294
$ws.Range("R6").Value = -730
$ws.Range("S6").Value = -1942
$ws.Range("T6").Value = 619
$ws.Range("U6").Value = 1675
$ws.Range("V6").Value = 26869
$ws.Range("W6").Value = 2.43
$ws.Range("X6").Value = -5.1
$ws.Range("Y6").Value = -14.25
$ws.Range("Z6").Value = -5.15
$ws.Range("AA6").Value = 181.97
$ws.Range("AB6").Value = 99.38
$ws.Range("AC6").Value = -3121
$ws.Range("AD6").Value = -2.32
$ws.Range("AE6").Value = 20372
$ws.Range("AF6").Value = 0.35
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 95432737

# Row 7
$ws.Range("D7").Value = 56873
$ws.Range("E7").Value = 2213
$ws.Range("G7").Value = -95
$ws.Range("H7").Value = -197
$ws.Range("I7").Value = -184
$ws.Range("K7").Value = 56343
$ws.Range("L7").Value = 36146
$ws.Range("M7").Value = 20197
$ws.Range("N7").Value = 19146
$ws.Range("P7").Value = 5891
$ws.Range("Q7").Value = 1457
$ws.Range("R7").Value = -1050
$ws.Range("S7").Value = -446
$ws.Range("T7").Value = 701
$ws.Range("U7").Value = 1464
$ws.Range("W7").Value = 3.89
$ws.Range("X7").Value = -0.35
$ws.Range("Y7").Value = -0.96
$ws.Range("Z7").Value = -0.35
$ws.Range("AA7").Value = 178.97
$ws.Range("AC7").Value = -193
$ws.Range("AD7").Value = -25.46
$ws.Range("AE7").Value = 20255
$ws.Range("AF7").Value = 0.24
$ws.Range("AG7").Value = 17
$ws.Range("AH7").Value = 0.34
$ws.Range("AI7").Value = -8.619999999999999

# Row 8
$ws.Range("D8").Value = 54491
$ws.Range("E8").Value = 1860
$ws.Range("G8").Value = 603
$ws.Range("H8").Value = 497
$ws.Range("I8").Value = 454
$ws.Range("K8").Value = 56150
$ws.Range("L8").Value = 35471
$ws.Range("M8").Value = 20679
$ws.Range("N8").Value = 19618
$ws.Range("P8").Value = 5891
$ws.Range("Q8").Value = 2719
$ws.Range("R8").Value = -1295
$ws.Range("S8").Value = -797
$ws.Range("T8").Value = 920
$ws.Range("U8").Value = 2264
$ws.Range("W8").Value = 3.41
$ws.Range("X8").Value = 0.91
$ws.Range("Y8").Value = 2.34
$ws.Range("Z8").Value = 0.88
$ws.Range("AA8").Value = 171.53
$ws.Range("AC8").Value = 476
$ws.Range("AD8").Value = 10.34
$ws.Range("AE8").Value = 20754
$ws.Range("AF8").Value = 0.24
$ws.Range("AG8").Value = 33
$ws.Range("AH8").Value = 0.68
$ws.Range("AI8").Value = 7

# Row 9
$ws.Range("D9").Value = 55264
$ws.Range("E9").Value = 2021
$ws.Range("G9").Value = 848
$ws.Range("H9").Value = 651
$ws.Range("I9").Value = 636
$ws.Range("K9").Value = 55770
$ws.Range("L9").Value = 34476
$ws.Range("M9").Value = 21294
$ws.Range("N9").Value = 20221
$ws.Range("P9").Value = 5891
$ws.Range("Q9").Value = 2175
$ws.Range("R9").Value = -1201
$ws.Range("S9").Value = -1224
$ws.Range("T9").Value = 920
$ws.Range("U9").Value = 2108
$ws.Range("W9").Value = 3.66
$ws.Range("X9").Value = 1.18
$ws.Range("Y9").Value = 3.19
$ws.Range("Z9").Value = 1.16
$ws.Range("AA9").Value = 161.91
$ws.Range("AC9").Value = 667
$ws.Range("AD9").Value = 7.38
$ws.Range("AE9").Value = 21392
$ws.Range("AF9").Value = 0.23
$ws.Range("AG9").Value = 33
$ws.Range("AH9").Value = 0.68
$ws.Range("AI9").Value = 5

# Cells removed entirely (no longer reported for this period)
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
